$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.459988594055176
$ws.Range("B1").Value = 3.622408866882324
$ws.Range("C1").Value = 6.173627376556396
$ws.Range("D1").Value = 1.505964040756226
$ws.Range("E1").Value = 0.8820995688438416
